# Raul's Log - add new AV-tech log entries (rows 1041-1049, 1054-1064)
# on the "Logs" worksheet, matching the upstream commit's row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Each record: Row, A=Task Type, B=Date (serial), C=Time, D=Building, E=Room, F=Comments (optional)
$rows = @(
    @{ Row=1041; A="Pickup PC";            B=42697; C="1730"; D="ATK"; E="005";  F="Return to DB 0003 and plug in for updates !!" }
    @{ Row=1042; A="Pickup Projector";     B=42697; C="1730"; D="ATK"; E="005";  F="Return to ATK 003C" }
    @{ Row=1043; A="Setup PC";             B=42697; C="1630"; D="HNE"; E="105";  F="Equipment from HNES 003" }
    @{ Row=1044; A="Pickup PC";            B=42697; C="1730"; D="HNE"; E="105";  F="Return equipment to HNES 003" }
    @{ Row=1045; A="Demo";                 B=42697; C="1745"; D="HNE"; E="036";  F=$null }
    @{ Row=1046; A="Demo";                 B=42697; C="1900"; D="SSB"; E="W133"; F=$null }
    @{ Row=1047; A="Pickup Mic";           B=42697; C="1715"; D="HNE"; E="402";  F="Return 4 IR mics with mixer and IR receivers to HNES 003" }
    @{ Row=1048; A="SCLD Student Event";   B=42697; C="1600"; D="FC";  E="109";  F="750624" }
    @{ Row=1049; A="SCLD Student Logout";  B=42697; C="1900"; D="FC";  E="109";  F="750624" }
    @{ Row=1054; A="Setup Mic";            B=42698; C="1800"; D="DB";  E="2027"; F="Neck mic and small PA from DB 0003" }
    @{ Row=1055; A="Pickup Mic";           B=42698; C="2100"; D="DB";  E="2027"; F="Return neck mic and small PA to DB 0003" }
    @{ Row=1056; A="Demo";                 B=42698; C="1900"; D="SSB"; E="N107"; F=$null }
    @{ Row=1057; A="Demo";                 B=42698; C="1900"; D="SSB"; E="S126"; F=$null }
    @{ Row=1058; A="Setup Mic";            B=42698; C="1630"; D="SSB"; E="W141"; F="4 desk mics from the rear booth - plug into floor box  / podium mic there -check"; Height=30 }
    @{ Row=1059; A="Operator";             B=42698; C="1630"; D="SSB"; E="W141"; F="Operate event from 4:30 to 7:30" }
    @{ Row=1060; A="Pickup Mic";           B=42698; C="1930"; D="SSB"; E="W141"; F="Return 4 desk mics with stands and cables to rear booth" }
    @{ Row=1061; A="SCLD Student Event";   B=42697; C="1800"; D="CLH"; E="K";    F="749658" }
    @{ Row=1062; A="SCLD Student Logout";  B=42697; C="2030"; D="CLH"; E="K";    F="749658" }
    @{ Row=1063; A="SCLD Student Event";   B=42697; C="1800"; D="ACW"; E="006";  F="751159" }
    @{ Row=1064; A="SCLD Student Logout";  B=42697; C="2030"; D="ACW"; E="006";  F="751159" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    if ($r.F -ne $null) {
        $ws.Cells.Item($row, 6).Value = $r.F
    }
    if ($r.Height -ne $null) {
        $ws.Rows.Item($row).RowHeight = $r.Height
    }
}

# Move the viewport/selection to mirror where the log was last edited.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1039
$ws.Range("F1068").Select() | Out-Null

Write-Host "Added 20 new log rows to 'Logs' sheet."
